# Predicted and plotted the air temperature distribution
# - Inlet/Outlet temps (G/H) recomputed using a more precise C->K offset
#   (273.15 instead of 273), which ripples into Air Density (J),
#   Mass Flow Rate (K) and Heat Transfer (BMT) (Q).
# - Relabelled "Bulk Mean Temperature" header -> "Bulk Mean Temperature Slope".
# - Shortened "Current" -> "Crnt" in the four test-case names.
# - Widened column A and column P to fit the new header/label text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / label text updates -------------------------------------------------
$ws.Range("P1").Value = "Bulk Mean Temperature Slope"

$ws.Range("A2").Value = "Test 1 - High MFR Low Crnt"
$ws.Range("A3").Value = "Test 2 - High MFR High Crnt"
$ws.Range("A4").Value = "Test 3 - Low MFR High Crnt"
$ws.Range("A5").Value = "Test 4 - Low MFR Low Crnt"

# --- Column width updates ---------------------------------------------------------
# (engine quantizes ColumnWidth to 1/6-character steps, so 26.8 is the closest
#  input that lands on the intended ~27.71-character rendered width)
$ws.Columns.Item(1).ColumnWidth = 26.8
$ws.Columns.Item(16).ColumnWidth = 26.8

# --- Recomputed data (Inlet Temp / Outlet Temp now use +273.15 K offset) ---------
# Row 2
$ws.Range("G2").Value = 310.65
$ws.Range("H2").Value = 325.15
$ws.Range("J2").Value = 1.198844968253475
$ws.Range("K2").Value = 0.07435630621823901
$ws.Range("Q2").Value = 1423.708655898139

# Row 3
$ws.Range("G3").Value = 312.65
$ws.Range("H3").Value = 340.15
$ws.Range("J3").Value = 1.190760517326973
$ws.Range("K3").Value = 0.07410517000335579
$ws.Range("Q3").Value = 3664.299599506274

# Row 4
$ws.Range("G4").Value = 313.95
$ws.Range("H4").Value = 347.15
$ws.Range("J4").Value = 1.157967055021532
$ws.Range("K4").Value = 0.05402724157953815
$ws.Range("Q4").Value = 3576.387522571998

# Row 5
$ws.Range("G5").Value = 314.32
$ws.Range("H5").Value = 332.15
$ws.Range("J5").Value = 1.15701727720054
$ws.Range("K5").Value = 0.05431106039320174
$ws.Range("Q5").Value = 1561.97602247048
